$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I width: 14 -> 10
$ws.Columns.Item(9).ColumnWidth = 10

# --- Sessions that moved from "Not Recorded" to "Recorded": restyle A:I to the same
#     highlight already used by recorded rows (copy format only, e.g. from row 2),
#     then refresh the Recorded-By / Students / Status values.
$styleSource = $ws.Range("A2:I2")
$flippedRows = 25, 44, 63, 172, 191, 210
foreach ($r in $flippedRows) {
    $styleSource.Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("G25").Value = "dnasr281@gmail.com"
$ws.Range("H25").Value = "22/31"
$ws.Range("I25").Value = "Recorded"

$ws.Range("G44").Value = "dnasr281@gmail.com"
$ws.Range("H44").Value = "13/19"
$ws.Range("I44").Value = "Recorded"

$ws.Range("G63").Value = "dnasr281@gmail.com"
$ws.Range("H63").Value = "19/21"
$ws.Range("I63").Value = "Recorded"

$ws.Range("G172").Value = "dnasr281@gmail.com"
$ws.Range("H172").Value = "24/27"
$ws.Range("I172").Value = "Recorded"

$ws.Range("G191").Value = "dnasr281@gmail.com"
$ws.Range("H191").Value = "22/29"
$ws.Range("I191").Value = "Recorded"

$ws.Range("G210").Value = "dnasr281@gmail.com"
$ws.Range("H210").Value = "22/29"
$ws.Range("I210").Value = "Recorded"

# --- Class Statistics summary (K/L column) - plain numbers
$ws.Range("L6").Value = 78
$ws.Range("L7").Value = 0

# --- Percentage cells are stored as literal text (e.g. "35.1%"), not numeric percentages,
#     in this workbook. Assigning a "nn.n%" string straight to a cell makes Excel coerce it
#     into a real percentage number, which would also mint a brand-new number-format style.
#     Route the text through a scratch cell that's pre-formatted as Text ("@") and bring the
#     value across with a values-only paste, which keeps the destination cell's original style
#     untouched and preserves the literal string.
$stage = $ws.Range("AA1")
$stage.NumberFormat = "@"

function Set-TextValue([string]$cellRef, [string]$val) {
    $stage.Value = $val
    $stage.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextValue "L9" "35.1%"
Set-TextValue "L10" "77.5%"

# --- Group Statistics table (M:S columns) for the groups whose recorded-session
#     counts moved because of the newly-recorded sessions above.
$newS = @{16="73.7%"; 17="58.6%"; 18="83.7%"; 24="69.8%"; 25="73.9%"; 26="69.5%"}
foreach ($r in 16, 17, 18, 24, 25, 26) {
    $ws.Range("O" + $r).Value = 7
    $ws.Range("P" + $r).Value = 0
    Set-TextValue ("R" + $r) "36.8%"
    Set-TextValue ("S" + $r) $newS[$r]
}

$stage.Clear()
$excel.CutCopyMode = $false
